# Reorganize the "Estados de Cuenta" (EC) rows: remove the previous period
# entries and add the new ones, regrouping each worker's 4 periods
# (2406, 2405, 2404, 2403) in descending order, per the updated source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico

$data = @(
    @{ Row = 16; Doc = "45578906";   Nombre = "DARLYS DEL SOCORRO MERCADO GARIZAO";  Periodo = "2406"; Mora = 52000 },
    @{ Row = 17; Doc = "45578906";   Nombre = "DARLYS DEL SOCORRO MERCADO GARIZAO";  Periodo = "2405"; Mora = 52000 },
    @{ Row = 18; Doc = "45578906";   Nombre = "DARLYS DEL SOCORRO MERCADO GARIZAO";  Periodo = "2404"; Mora = 52000 },
    @{ Row = 19; Doc = "45578906";   Nombre = "DARLYS DEL SOCORRO MERCADO GARIZAO";  Periodo = "2403"; Mora = 52000 },
    @{ Row = 20; Doc = "22422803";   Nombre = "CECILIA DE LA CRUZ TORRES HERNANDEZ"; Periodo = "2406"; Mora = 52000 },
    @{ Row = 21; Doc = "22422803";   Nombre = "CECILIA DE LA CRUZ TORRES HERNANDEZ"; Periodo = "2405"; Mora = 52000 },
    @{ Row = 22; Doc = "22422803";   Nombre = "CECILIA DE LA CRUZ TORRES HERNANDEZ"; Periodo = "2404"; Mora = 10400 },
    @{ Row = 23; Doc = "32939137";   Nombre = "ELIANA VERONICA BATISTA MARTELO";     Periodo = "2406"; Mora = 52000 },
    @{ Row = 24; Doc = "32939137";   Nombre = "ELIANA VERONICA BATISTA MARTELO";     Periodo = "2405"; Mora = 52000 },
    @{ Row = 25; Doc = "32939137";   Nombre = "ELIANA VERONICA BATISTA MARTELO";     Periodo = "2404"; Mora = 10400 },
    @{ Row = 26; Doc = "1100339103"; Nombre = "JESUS MANUEL GONZALEZ CONDE";         Periodo = "2406"; Mora = 31200 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Mora
}
